$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-17 (16 rows x 20 columns, A..T)
$data = @(
    @("ECs", "Sema3c", "Nrp2", "ECs", [double]3, [double]1, [double]"0.604961", [double]"1.814883", [double]"0.0160517527720356", [double]"0.0160517527720356", [double]3, [double]1, [double]"57.65261933333333", [double]"172.957858", [double]"0.6817060950001529", [double]"0.6817060950001529", [double]"34.87758624451266", [double]"313.898276200614", [double]"0.01094257770013227", [double]"0.01094257770013227"),
    @("ECs", "Sema3c", "Nrp2", "FAPs", [double]3, [double]1, [double]"0.604961", [double]"1.814883", [double]"0.0160517527720356", [double]"0.0160517527720356", [double]3, [double]1, [double]"8.622273333333332", [double]"25.86682", [double]"0.1019529789289588", [double]"0.1019529789289588", [double]"5.216139098006666", [double]"46.94525188205999", [double]"0.001636524012140202", [double]"0.001636524012140202"),
    @("ECs", "Sema3c", "Nrp2", "MuSCs", [double]3, [double]1, [double]"0.604961", [double]"1.814883", [double]"0.0160517527720356", [double]"0.0160517527720356", [double]3, [double]1, [double]"2.790736", [double]"8.372208", [double]"0.0329987043561157", [double]"0.0329987043561157", [double]"1.688286441296", [double]"15.194577971664", [double]"0.0005296870441218635", [double]"0.0005296870441218635"),
    @("ECs", "Sema3c", "Nrp2", "Resolving-Mac", [double]3, [double]1, [double]"0.604961", [double]"1.814883", [double]"0.0160517527720356", [double]"0.0160517527720356", [double]3, [double]1, [double]"15.50544933333333", [double]"46.516348", [double]"0.1833422217147727", [double]"0.1833422217147727", [double]"9.380192134142666", [double]"84.421729207284", [double]"0.002942964015641268", [double]"0.002942964015641269"),
    @("FAPs", "Sema3c", "Nrp2", "ECs", [double]3, [double]1, [double]"35.10542566666667", [double]"105.316277", [double]"0.9314709770686151", [double]"0.9314709770686151", [double]3, [double]1, [double]"57.65261933333333", [double]"172.957858", [double]"0.6817060950001529", [double]"0.6817060950001529", [double]"2023.919742494963", [double]"18215.27768245466", [double]"0.6349894423834226", [double]"0.6349894423834226"),
    @("FAPs", "Sema3c", "Nrp2", "FAPs", [double]3, [double]1, [double]"35.10542566666667", [double]"105.316277", [double]"0.9314709770686151", [double]"0.9314709770686151", [double]3, [double]1, [double]"8.622273333333332", [double]"25.86682", [double]"0.1019529789289588", [double]"0.1019529789289588", [double]"302.6885755810155", [double]"2724.19718022914", [double]"0.09496624089801321", [double]"0.09496624089801321"),
    @("FAPs", "Sema3c", "Nrp2", "MuSCs", [double]3, [double]1, [double]"35.10542566666667", [double]"105.316277", [double]"0.9314709770686151", [double]"0.9314709770686151", [double]3, [double]1, [double]"2.790736", [double]"8.372208", [double]"0.0329987043561157", [double]"0.0329987043561157", [double]"97.96997520329069", [double]"881.7297768296161", [double]"0.03073733538858946", [double]"0.03073733538858945"),
    @("FAPs", "Sema3c", "Nrp2", "Resolving-Mac", [double]3, [double]1, [double]"35.10542566666667", [double]"105.316277", [double]"0.9314709770686151", [double]"0.9314709770686151", [double]3, [double]1, [double]"15.50544933333333", [double]"46.516348", [double]"0.1833422217147727", [double]"0.1833422217147727", [double]"544.3253989995995", [double]"4898.928590996396", [double]"0.17077795839859", [double]"0.17077795839859"),
    @("MuSCs", "Sema3c", "Nrp2", "ECs", [double]3, [double]1, [double]"1.958375", [double]"5.875125", [double]"0.05196260806057782", [double]"0.05196260806057783", [double]3, [double]1, [double]"57.65261933333333", [double]"172.957858", [double]"0.6817060950001529", [double]"0.6817060950001529", [double]"112.9054483869166", [double]"1016.14903548225", [double]"0.03542322662699997", [double]"0.03542322662699998"),
    @("MuSCs", "Sema3c", "Nrp2", "FAPs", [double]3, [double]1, [double]"1.958375", [double]"5.875125", [double]"0.05196260806057782", [double]"0.05196260806057783", [double]3, [double]1, [double]"8.622273333333332", [double]"25.86682", [double]"0.1019529789289588", [double]"0.1019529789289588", [double]"16.88564453916667", [double]"151.9708008525", [double]"0.005297742684693837", [double]"0.005297742684693837"),
    @("MuSCs", "Sema3c", "Nrp2", "MuSCs", [double]3, [double]1, [double]"1.958375", [double]"5.875125", [double]"0.05196260806057782", [double]"0.05196260806057783", [double]3, [double]1, [double]"2.790736", [double]"8.372208", [double]"0.0329987043561157", [double]"0.0329987043561157", [double]"5.465307614", [double]"49.187768526", [double]"0.001714698740963722", [double]"0.001714698740963722"),
    @("MuSCs", "Sema3c", "Nrp2", "Resolving-Mac", [double]3, [double]1, [double]"1.958375", [double]"5.875125", [double]"0.05196260806057782", [double]"0.05196260806057783", [double]3, [double]1, [double]"15.50544933333333", [double]"46.516348", [double]"0.1833422217147727", [double]"0.1833422217147727", [double]"30.36548433816667", [double]"273.2893590435", [double]"0.009526940007920292", [double]"0.009526940007920294"),
    @("Resolving-Mac", "Sema3c", "Nrp2", "ECs", [double]1, [double]"0.3333333333333333", [double]"0.01939666666666667", [double]"0.05819", [double]"0.0005146620987715195", [double]"0.0005146620987715195", [double]3, [double]1, [double]"57.65261933333333", [double]"172.957858", [double]"0.6817060950001529", [double]"0.6817060950001529", [double]"1.118268639668889", [double]"10.06441775702", [double]"0.0003508482895981156", [double]"0.0003508482895981156"),
    @("Resolving-Mac", "Sema3c", "Nrp2", "FAPs", [double]1, [double]"0.3333333333333333", [double]"0.01939666666666667", [double]"0.05819", [double]"0.0005146620987715195", [double]"0.0005146620987715195", [double]3, [double]1, [double]"8.622273333333332", [double]"25.86682", [double]"0.1019529789289588", [double]"0.1019529789289588", [double]"0.1672433617555555", [double]"1.5051902558", [double]"5.247133411158646e-05", [double]"5.247133411158646e-05"),
    @("Resolving-Mac", "Sema3c", "Nrp2", "MuSCs", [double]1, [double]"0.3333333333333333", [double]"0.01939666666666667", [double]"0.05819", [double]"0.0005146620987715195", [double]"0.0005146620987715195", [double]3, [double]1, [double]"2.790736", [double]"8.372208", [double]"0.0329987043561157", [double]"0.0329987043561157", [double]"0.05413097594666667", [double]"0.48717878352", [double]"1.698318244065939e-05", [double]"1.698318244065939e-05"),
    @("Resolving-Mac", "Sema3c", "Nrp2", "Resolving-Mac", [double]1, [double]"0.3333333333333333", [double]"0.01939666666666667", [double]"0.05819", [double]"0.0005146620987715195", [double]"0.0005146620987715195", [double]3, [double]1, [double]"15.50544933333333", [double]"46.516348", [double]"0.1833422217147727", [double]"0.1833422217147727", [double]"0.3007540322355555", [double]"2.70678629012", [double]"9.435929262115818e-05", [double]"9.435929262115818e-05"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Output "Done"